$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.060.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.420.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "410.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.00%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "44.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +34.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.77%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.963.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.437.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.51%  "
$ws.Range("E19").Value = "  +6.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.114.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "501.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +58.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "94.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.19%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.89%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.78%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0512"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +8.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "144.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.23%  "
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("E48").Value = "  +5.49%  "
$ws.Range("E49").Value = "  +15.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +36.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.98%  "
